# ---------------------------------------------------------------------------
# Commit: "feat: add 2022-Q3 data"
#
# 1) Insert a new 2022-Q3 row at the top of the "总计" (summary) sheet,
#    pushing existing quarters down by one row, and append the final
#    "2020-Q4" row that falls off the bottom of the old table.
# 2) Insert a brand-new "2022-Q3" worksheet (cloned from the structurally
#    identical "2022-Q2" sheet so headers/styles match exactly) right after
#    "总计", populated with the quarter's fund holdings table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Update "总计" summary sheet -----------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Insert a fresh row 2 for the new quarter; shift 2022-Q2..2020-Q4 down one.
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Give the new row2's A cell the same style used by the rest of column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("A9").PasteSpecial(-4122)

# New 2022-Q3 summary row.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 13
$summary.Range("D2").Value = 2.3

# Re-sequence the running index column A2:A9 = 0..7.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# The old row 8 (2020-Q4) now lives at row 9.
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 6
$summary.Range("D9").Value = 0.8100000000000001

# --- 2. Add the "2022-Q3" fund-holdings worksheet --------------------------
# Clone "2022-Q2" (same column layout/header/styles) and drop it in right
# after "总计", then rename and overwrite its data with the 2022-Q3 figures.
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$afterSheet = $wb.Worksheets.Item("总计")
$templateSheet.Copy($null, $afterSheet)
$ws = $wb.Worksheets.Item("2022-Q2 (2)")
$ws.Name = "2022-Q3"

# Template sheet had 16 data rows; the 2022-Q3 table only has 14, trim it.
$ws.Range("A15:H16").EntireRow.Delete()

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'100032"
$ws.Cells.Item(2,3).Value = "富国中证红利指数增强A"
$ws.Cells.Item(2,4).Value = "'54.79"
$ws.Cells.Item(2,5).Value = "'91.26"
$ws.Cells.Item(2,6).Value = "'1.66"
$ws.Cells.Item(2,7).Value = "'0.9095"
$ws.Cells.Item(2,8).Value = 9
$ws.Range("B2:G2").ClearFormats()

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'090010"
$ws.Cells.Item(3,3).Value = "大成中证红利指数A"
$ws.Cells.Item(3,4).Value = "'26.31"
$ws.Cells.Item(3,5).Value = "'93.17"
$ws.Cells.Item(3,6).Value = "'1.75"
$ws.Cells.Item(3,7).Value = "'0.4604"
$ws.Cells.Item(3,8).Value = 4
$ws.Range("B3:G3").ClearFormats()

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'515180"
$ws.Cells.Item(4,3).Value = "易方达中证红利ETF"
$ws.Cells.Item(4,4).Value = "'16.81"
$ws.Cells.Item(4,5).Value = "'99.41"
$ws.Cells.Item(4,6).Value = "'1.87"
$ws.Cells.Item(4,7).Value = "'0.3143"
$ws.Cells.Item(4,8).Value = 4
$ws.Range("B4:G4").ClearFormats()

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'515080"
$ws.Cells.Item(5,3).Value = "招商中证红利ETF"
$ws.Cells.Item(5,4).Value = "'8.33"
$ws.Cells.Item(5,5).Value = "'99.62"
$ws.Cells.Item(5,6).Value = "'1.87"
$ws.Cells.Item(5,7).Value = "'0.1558"
$ws.Cells.Item(5,8).Value = 4
$ws.Range("B5:G5").ClearFormats()

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'008682"
$ws.Cells.Item(6,3).Value = "富国中证红利指数增强C"
$ws.Cells.Item(6,4).Value = "'8.27"
$ws.Cells.Item(6,5).Value = "'91.26"
$ws.Cells.Item(6,6).Value = "'1.66"
$ws.Cells.Item(6,7).Value = "'0.1373"
$ws.Cells.Item(6,8).Value = 9
$ws.Range("B6:G6").ClearFormats()

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'501029"
$ws.Cells.Item(7,3).Value = "华宝标普中国A股红利机会指数（LOF）A"
$ws.Cells.Item(7,4).Value = "'8.11"
$ws.Cells.Item(7,5).Value = "'94.26"
$ws.Cells.Item(7,6).Value = "'1.67"
$ws.Cells.Item(7,7).Value = "'0.1354"
$ws.Cells.Item(7,8).Value = 4
$ws.Range("B7:G7").ClearFormats()

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'005125"
$ws.Cells.Item(8,3).Value = "华宝标普中国A股红利机会指数C"
$ws.Cells.Item(8,4).Value = "'3.38"
$ws.Cells.Item(8,5).Value = "'94.26"
$ws.Cells.Item(8,6).Value = "'1.67"
$ws.Cells.Item(8,7).Value = "'0.0564"
$ws.Cells.Item(8,8).Value = 4
$ws.Range("B8:G8").ClearFormats()

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'515890"
$ws.Cells.Item(9,3).Value = "博时中证红利ETF"
$ws.Cells.Item(9,4).Value = "'2.92"
$ws.Cells.Item(9,5).Value = "'98.00"
$ws.Cells.Item(9,6).Value = "'1.83"
$ws.Cells.Item(9,7).Value = "'0.0534"
$ws.Cells.Item(9,8).Value = 4
$ws.Range("B9:G9").ClearFormats()

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'007801"
$ws.Cells.Item(10,3).Value = "大成中证红利指数C"
$ws.Cells.Item(10,4).Value = "'2.72"
$ws.Cells.Item(10,5).Value = "'93.17"
$ws.Cells.Item(10,6).Value = "'1.75"
$ws.Cells.Item(10,7).Value = "'0.0476"
$ws.Cells.Item(10,8).Value = 4
$ws.Range("B10:G10").ClearFormats()

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'161907"
$ws.Cells.Item(11,3).Value = "万家中证红利指数（LOF）A"
$ws.Cells.Item(11,4).Value = "'1.30"
$ws.Cells.Item(11,5).Value = "'94.17"
$ws.Cells.Item(11,6).Value = "'1.82"
$ws.Cells.Item(11,7).Value = "'0.0237"
$ws.Cells.Item(11,8).Value = 4
$ws.Range("B11:G11").ClearFormats()

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "'501307"
$ws.Cells.Item(12,3).Value = "银河中证沪港深高股息指数（LOF）A"
$ws.Cells.Item(12,4).Value = "'0.15"
$ws.Cells.Item(12,5).Value = "'90.33"
$ws.Cells.Item(12,6).Value = "'1.50"
$ws.Cells.Item(12,7).Value = "'0.0022"
$ws.Cells.Item(12,8).Value = 5
$ws.Range("B12:G12").ClearFormats()

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "'015558"
$ws.Cells.Item(13,3).Value = "万家中证红利指数（LOF）C"
$ws.Cells.Item(13,4).Value = "'0.01"
$ws.Cells.Item(13,5).Value = "'94.17"
$ws.Cells.Item(13,6).Value = "'1.82"
$ws.Cells.Item(13,7).Value = "'0.0002"
$ws.Cells.Item(13,8).Value = 4
$ws.Range("B13:G13").ClearFormats()

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "'501308"
$ws.Cells.Item(14,3).Value = "银河中证沪港深高股息指数（LOF）C"
$ws.Cells.Item(14,4).Value = "'0.01"
$ws.Cells.Item(14,5).Value = "'90.33"
$ws.Cells.Item(14,6).Value = "'1.50"
$ws.Cells.Item(14,7).Value = "'0.0002"
$ws.Cells.Item(14,8).Value = 5
$ws.Range("B14:G14").ClearFormats()

Write-Output "2022-Q3 sheet inserted and 总计 summary updated"
